$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Initial Position Single Thread")
$ws2 = $wb.Worksheets.Item("vs other Enignes")

# ---------------------------------------------------------------------------
# 1) Build new rows 109-111 on "Initial Position Single Thread" by cloning
#    the formatting of the analogous existing rows (105-107), then filling
#    in the new values / formulas (mirrors a new "commit" block being added
#    under row 105-107's block).
# ---------------------------------------------------------------------------

# Row 109 formatting <- Row 105 (full block start row, columns A:N)
$ws.Range("A105:N105").Copy()
$ws.Range("A109").PasteSpecial(-4122)

# Row 110 formatting <- Row 106 (columns C:N)
$ws.Range("C106:N106").Copy()
$ws.Range("C110").PasteSpecial(-4122)

# Row 111 formatting <- Row 107 (columns I:N)
$ws.Range("I107:N107").Copy()
$ws.Range("I111").PasteSpecial(-4122)

# Row 109 values / formulas
$ws.Range("A109").Value = 46065
$ws.Range("C109").Value = 4
$ws.Range("D109").Value = 206603
$ws.Range("E109").Value = 147
$ws.Range("F109").Formula = "=D109/E109*1000"
$ws.Range("G109").Formula = "=(E100-E109)/E100"
$ws.Range("H109").Formula = "=(F109-80000000)/80000000"
$ws.Range("I109").Value = 4
$ws.Range("J109").Value = 197281
$ws.Range("K109").Value = 6
$ws.Range("L109").Formula = "=J109/K109*1000"
$ws.Range("M109").Formula = "=(K100-K109)/K100"
$ws.Range("N109").Formula = "=(L109-80000000)/80000000"

# Row 110 values / formulas
$ws.Range("C110").Value = 5
$ws.Range("D110").Value = 5072212
$ws.Range("E110").Value = 3204
$ws.Range("F110").Formula = "=D110/E110*1000"
$ws.Range("G110").Formula = "=(E101-E110)/E101"
$ws.Range("H110").Formula = "=(F110-80000000)/80000000"
$ws.Range("I110").Value = 5
$ws.Range("J110").Value = 4880523
$ws.Range("K110").Value = 115
$ws.Range("L110").Formula = "=J110/K110*1000"
$ws.Range("M110").Formula = "=(K101-K110)/K101"
$ws.Range("N110").Formula = "=(L110-80000000)/80000000"

# Row 111 values / formulas
$ws.Range("I111").Value = 6
$ws.Range("J111").Value = 119060324
$ws.Range("K111").Value = 2536
$ws.Range("L111").Formula = "=J111/K111*1000"
$ws.Range("M111").Formula = "=(K102-K111)/K102"
$ws.Range("N111").Formula = "=(L111-80000000)/80000000"

# ---------------------------------------------------------------------------
# 2) "note" column (P): the newest commit note moves from row 105 to the new
#    row 109, and row 105's note cell reverts to the regular (non-highlight)
#    style used by older rows such as P100.
# ---------------------------------------------------------------------------

# Give P109 the highlighted "latest" style that P105 currently has ...
$ws.Range("P105").Copy()
$ws.Range("P109").PasteSpecial(-4122)
$ws.Range("P109").Value = "Improved compilation flags"

# ... then restyle P105 back to the plain note style (same as P100).
$ws.Range("P100").Copy()
$ws.Range("P105").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Update the view so the new rows are visible / selected.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 73
$ws.Range("E110").Select()

# ---------------------------------------------------------------------------
# 4) "vs other Enignes" sheet: B1/C1/D1 are untouched content-wise (still
#    "stockfish livello 1" / "vinto" / "https://lichess.org/") - re-assert
#    them defensively so the sheet keeps referencing those exact strings.
# ---------------------------------------------------------------------------
$ws2.Range("B1").Value = "stockfish livello 1"
$ws2.Range("C1").Value = "vinto"
$ws2.Range("D1").Value = "https://lichess.org/"
